# Update the two-digit multiplication problems/answers in the table.
$d = $word.ActiveDocument

$replacements = @(
    @{ old = "92×52=4784"; new = "85×68=5780" },
    @{ old = "11×46=506";  new = "36×50=1800" },
    @{ old = "64×69=4416"; new = "64×20=1280" },
    @{ old = "47×26=1222"; new = "14×79=1106" },
    @{ old = "54×87=4698"; new = "44×79=3476" },
    @{ old = "15×14=210";  new = "88×82=7216" },
    @{ old = "48×12=576";  new = "71×69=4899" },
    @{ old = "56×83=4648"; new = "85×28=2380" },
    @{ old = "66×19=1254"; new = "96×72=6912" },
    @{ old = "90×13=1170"; new = "79×95=7505" },
    @{ old = "92×25=2300"; new = "47×11=517"  },
    @{ old = "95×89=8455"; new = "52×93=4836" },
    @{ old = "19×35=665";  new = "80×89=7120" },
    @{ old = "25×84=2100"; new = "76×56=4256" },
    @{ old = "82×23=1886"; new = "24×19=456"  },
    @{ old = "46×49=2254"; new = "42×20=840"  },
    @{ old = "84×77=6468"; new = "35×94=3290" },
    @{ old = "89×33=2937"; new = "41×92=3772" },
    @{ old = "60×39=2340"; new = "35×55=1925" },
    @{ old = "17×56=952";  new = "72×35=2520" },
    @{ old = "62×90=5580"; new = "74×14=1036" },
    @{ old = "95×28=2660"; new = "48×50=2400" },
    @{ old = "36×37=1332"; new = "98×47=4606" },
    @{ old = "35×98=3430"; new = "77×64=4928" },
    @{ old = "27×76=2052"; new = "34×21=714"  }
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
